$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = -602.4
$ws.Range("B3").Value = -965.7
$ws.Range("C3").Value = -654.2
$ws.Range("C4").Value = -572.9
$ws.Range("C9").Value = -543.5
$ws.Range("C23").Value = -416.5
